# Update horarios-141-completo.xlsx: append newly scraped rows (31/12/2025, scrape
# run finishing 17:50:18) to the three sheets, and refresh the "last updated" /
# "total rows" header cells on each sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "LP1912" (sheet1): columns A Fecha(blank)/B Hora_Scrap/C Hora_Llegada/
# D Linea/E Minutos/F Parada/G Fecha -- append rows 1164..1181
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: 31/12/2025 17:50:18"
$ws1.Range("A3").Value = "Total filas: 1180"

$sheet1Rows = @(
    @(1164, "17:50:07", "17:55", "10_OLMOS", 5, "LP1912", "31/12/2025"),
    @(1165, "17:50:07", "18:00", "16_SANTA ANA", 10, "LP1912", "31/12/2025"),
    @(1166, "17:50:07", "18:04", "23_HERNANDEZ", 14, "LP1912", "31/12/2025"),
    @(1167, "17:50:07", "18:05", "14_ABASTO", 15, "LP1912", "31/12/2025"),
    @(1168, "17:50:07", "18:22", "16_SANTA ANA", 32, "LP1912", "31/12/2025"),
    @(1169, "17:50:07", "18:25", "11_ETCHEVERRY", 35, "LP1912", "31/12/2025"),
    @(1170, "17:50:07", "18:28", "15_ABASTO", 38, "LP1912", "31/12/2025"),
    @(1171, "17:50:07", "18:34", "14X44_ABASTO", 44, "LP1912", "31/12/2025"),
    @(1172, "17:50:07", "18:34", "23_HERNANDEZ", 44, "LP1912", "31/12/2025"),
    @(1173, "17:50:07", "18:42", "16_P MOR-SANTA ANA", 52, "LP1912", "31/12/2025"),
    @(1174, "17:50:07", "18:45", "14_ABASTO", 55, "LP1912", "31/12/2025"),
    @(1175, "17:50:07", "18:52", "15_ABASTO", 62, "LP1912", "31/12/2025"),
    @(1176, "17:50:07", "18:59", "10_OLMOS", 69, "LP1912", "31/12/2025"),
    @(1177, "17:50:07", "19:00", "23_HERNANDEZ", 70, "LP1912", "31/12/2025"),
    @(1178, "17:50:07", "19:02", "17_ROMERO", 72, "LP1912", "31/12/2025"),
    @(1179, "17:50:07", "19:07", "14_ABASTO", 77, "LP1912", "31/12/2025"),
    @(1180, "17:50:07", "19:12", "81_EL PELIGRO", 82, "LP1912", "31/12/2025"),
    @(1181, "17:50:07", "19:21", "215C_EL PATO", 91, "LP1912", "31/12/2025")
)

foreach ($row in $sheet1Rows) {
    $r = $row[0]
    $ws1.Cells.Item($r, 2).Value = $row[1]
    $ws1.Cells.Item($r, 3).Value = $row[2]
    $ws1.Cells.Item($r, 4).Value = $row[3]
    $ws1.Cells.Item($r, 5).Value = $row[4]
    $ws1.Cells.Item($r, 6).Value = $row[5]
    $ws1.Cells.Item($r, 7).Value = $row[6]
}

# ---------------------------------------------------------------------------
# Sheet "LP1912-215" (sheet2): columns A(blank)/B Fecha/C Hora_Scrap/
# D Hora_Llegada/E Linea/F Minutos/G Parada -- append row 79
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = "Última actualización: 31/12/2025 17:50:18"
$ws2.Range("A3").Value = "Total filas: 78"

$ws2.Cells.Item(79, 2).Value = "31/12/2025"
$ws2.Cells.Item(79, 3).Value = "17:50:07"
$ws2.Cells.Item(79, 4).Value = "19:21"
$ws2.Cells.Item(79, 5).Value = "215C_EL PATO"
$ws2.Cells.Item(79, 6).Value = 91
$ws2.Cells.Item(79, 7).Value = "LP1912"

# ---------------------------------------------------------------------------
# Sheet "6203-6173" (sheet3): columns A(blank)/B Fecha/C Hora_Scrap/
# D Hora_Llegada/E Linea/F Minutos/G Parada -- append rows 141..142
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = "Última actualización: 31/12/2025 17:50:18"
$ws3.Range("A3").Value = "Total filas: 141"

$ws3.Cells.Item(141, 2).Value = "31/12/2025"
$ws3.Cells.Item(141, 3).Value = "17:50:13"
$ws3.Cells.Item(141, 4).Value = "18:22"
$ws3.Cells.Item(141, 5).Value = "215C_LA PLATA"
$ws3.Cells.Item(141, 6).Value = 32
$ws3.Cells.Item(141, 7).Value = "L6203"

$ws3.Cells.Item(142, 2).Value = "31/12/2025"
$ws3.Cells.Item(142, 3).Value = "17:50:18"
$ws3.Cells.Item(142, 4).Value = "19:11"
$ws3.Cells.Item(142, 5).Value = "215B_LP-P MOR-1 Y 57"
$ws3.Cells.Item(142, 6).Value = 81
$ws3.Cells.Item(142, 7).Value = "L6173"
